$wb = $excel.ActiveWorkbook

# --- Sheet 1: "peliculas o documentales" ---
$ws1 = $wb.Worksheets.Item(1)

# Row 13: F13 currently holds the text "2019" (it was typed/stored as a shared
# string). Make it a genuine number while preserving the cell's existing
# style (s="1", text-format @). We briefly switch the cell to the default
# "Normal" style (General format) so the numeric literal isn't re-coerced to
# text, write the number, then restore the original "@" number format.
$ws1.Range("F13").Style = "Normal"
$ws1.Range("F13").Value = 2019
$ws1.Range("F13").NumberFormat = "@"

# New row 14 - "Guerrilla del Oro" (repeat entry)
$ws1.Range("A14").Value = "Guerrilla del Oro"
$ws1.Range("B14").Value = "NATIONAL GEOGRAPHIC"
$ws1.Range("C14").Value = "Premium"
$ws1.Range("D14").Value = "Documental"
$ws1.Range("E14").Value = "HD"
$ws1.Range("F14").Style = "Normal"
$ws1.Range("F14").Value = 2019
$ws1.Range("F14").NumberFormat = "@"

# New row 15 - "Sumergidos" (repeat entry)
$ws1.Range("A15").Value = "Sumergidos"
$ws1.Range("C15").Value = "Arriendo"
$ws1.Range("D15").Value = "Acción"
$ws1.Range("E15").Value = "HD"
$ws1.Range("F15").Style = "Normal"
$ws1.Range("F15").Value = 2019
$ws1.Range("F15").NumberFormat = "@"
$ws1.Range("G15").Value = "$3.490"

# --- Sheet 2: "series" ---
$ws2 = $wb.Worksheets.Item(2)

# New row 7 - "The Outsider - El Visitante"
$ws2.Range("A7").Value = "The Outsider - El Visitante"
$ws2.Range("C7").Value = "Gratis"
$ws2.Range("D7").Value = "Acción"
$ws2.Range("E7").Value = "HD"
$ws2.Range("F7").Value = "2020"
$ws2.Range("G7").Value = "Temporada 01"
$ws2.Range("H7").Value = "10"
